# "Loan RBI, Variable Instalments"
#
# The "Repayment schedule" sheet gains a new (blank) column between the
# existing "Late" / "heading" / "Outstanding" columns and the rest of the
# table (an inserted column N, shifting the old N/O/P -> O/P/Q). The
# "Edit Repayment Schedule" sheet's selection moves, and the
# "Repayment schedule" sheet becomes the active tab with a new selection.

$wb = $excel.ActiveWorkbook

# --- Repayment schedule: insert a new blank column before column N -------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("N:N").Insert()
# Column N keeps a fixed (non bestFit) width of 11, matching column M.
$wsSchedule.Columns("N:N").ColumnWidth = 10.166666666666666

# --- Edit Repayment Schedule: move the selection to B5 -------------------
$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsEdit.Range("B5").Select()

# --- Make "Repayment schedule" the active sheet/tab with selection K16 ---
$wsSchedule.Activate()
$wsSchedule.Range("K16").Select()
